# Add a new "2020" column (Q) to the 1.a.2 indicator table, mirroring the
# formatting already used by the neighbouring "2019" column (P) for each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (year header): 2020, formatted like P4 (year 2019 header).
$ws.Range("P4").Copy()
$ws.Range("Q4").PasteSpecial(-4122)
$ws.Range("Q4").Value = 2020

# Row 5 (Education): 53.2, formatted like P5.
$ws.Range("P5").Copy()
$ws.Range("Q5").PasteSpecial(-4122)
$ws.Range("Q5").Value = 53.2

# Row 6 (Health): 23.2, formatted like P6.
$ws.Range("P6").Copy()
$ws.Range("Q6").PasteSpecial(-4122)
$ws.Range("Q6").Value = 23.2

# Row 7 (Social protection): 10, formatted like P7 but with a "0.0" number
# format applied (so it reads consistently with the rest of the row).
$ws.Range("P7").Copy()
$ws.Range("Q7").PasteSpecial(-4122)
$ws.Range("Q7").NumberFormat = "0.0"
$ws.Range("Q7").Value = 10

# Row 8 (bottom, bordered row): 20, formatted like P8 but with the "0.0"
# number format applied.
$ws.Range("P8").Copy()
$ws.Range("Q8").PasteSpecial(-4122)
$ws.Range("Q8").NumberFormat = "0.0"
$ws.Range("Q8").Value = 20

# Match the saved selection state (cell P9 selected).
$null = $ws.Range("P9").Select()
